$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Offerings")
$ws.Select()
$ws.Range("B14").Select()
